$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1. New data for rows 132-143 (OBI / IAO / UO terms added per BBRB tracker
#    issue 25), plus the "yes" flag added at E16.
# ---------------------------------------------------------------------------

$ws.Range("E16").Value = "yes"

$rows = @(
    @{r=132; a="http://purl.obolibrary.org/obo/OBI_0001064"; b="automatic tissue processor"},
    @{r=133; a="http://purl.obolibrary.org/obo/OBI_0400168"; b="microtome"},
    @{r=134; a="http://purl.obolibrary.org/obo/OBI_0000852"; b="record of missing knowledge"},
    @{r=135; a="http://purl.obolibrary.org/obo/IAO_0000129"; b="version number"},
    @{r=136; a="http://purl.obolibrary.org/obo/IAO_0000329"; b="running title"},
    @{r=137; a="http://purl.obolibrary.org/obo/IAO_0000305"; b="document title"},
    @{r=138; a="http://purl.obolibrary.org/obo/IAO_0000414"; b="mass measurement datum"},
    @{r=139; a="http://purl.obolibrary.org/obo/UO_0000005";  b="temperature unit"},
    @{r=140; a="http://purl.obolibrary.org/obo/UO_0000027";  b="degree celsius"},
    @{r=141; a="http://purl.obolibrary.org/obo/UO_0000196";  b="pH"},
    @{r=142; a="http://purl.obolibrary.org/obo/OBI_0000835"; b="manufacturer"},
    @{r=143; a="http://purl.obolibrary.org/obo/IAO_0000132"; b="lot number"}
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Range("A$r").Value = $row.a
    $ws.Range("B$r").Value = $row.b
    $ws.Range("C$r").Value = "y"
}

# ---------------------------------------------------------------------------
# 1b. Hyperlinks for the new OBO term references in column A. These must be
#     added BEFORE the custom font formatting below: Hyperlinks.Add() resets
#     the target cell's font back to the workbook's default "Hyperlink"
#     cell style, which would otherwise clobber our explicit font choices.
# ---------------------------------------------------------------------------

foreach ($row in $rows) {
    $r = $row.r
    $ws.Hyperlinks.Add($ws.Range("A$r"), $row.a) | Out-Null
}

# ---------------------------------------------------------------------------
# 2. Formatting templates.
#    Column A (rows 132-138 and 139-143) = manual hyperlink-blue look.
#    Column B uses two different looks depending on the row block.
# ---------------------------------------------------------------------------

# Template for A132:A138  -> underline, 10pt blue Arial
$tA1 = $ws.Range("Z300")
$tA1.Font.Underline = 2
$tA1.Font.Size = 10
$tA1.Font.Color = 16711680
$tA1.Font.Name = "Arial"

# Template for A139:A143 -> underline, 11pt black Calibri
$tA2 = $ws.Range("Z301")
$tA2.Font.Underline = 2
$tA2.Font.Size = 11
$tA2.Font.Color = 0
$tA2.Font.Name = "Calibri"

# Template for B132:B134 -> 10pt Arial, wrap text
$tB1 = $ws.Range("Z302")
$tB1.Font.Size = 10
$tB1.Font.Name = "Arial"
$tB1.WrapText = $true

# Template for B135:B138 -> 11pt black Calibri, wrap text
$tB2 = $ws.Range("Z303")
$tB2.Font.Size = 11
$tB2.Font.Color = 0
$tB2.Font.Name = "Calibri"
$tB2.WrapText = $true

# Template for B139:B145 -> 10pt Arial, no wrap
$tB3 = $ws.Range("Z304")
$tB3.Font.Size = 10
$tB3.Font.Name = "Arial"

# Propagate formats without creating duplicate style/font table entries.
$tA1.Copy()
$ws.Range("A132:A138").PasteSpecial(-4122)

$tA2.Copy()
$ws.Range("A139:A143").PasteSpecial(-4122)

$tB1.Copy()
$ws.Range("B132:B134").PasteSpecial(-4122)

$tB2.Copy()
$ws.Range("B135:B138").PasteSpecial(-4122)

$tB3.Copy()
$ws.Range("B139:B145").PasteSpecial(-4122)

# Clear the scratch template cells (keep everything else untouched). Only
# clears contents/format of the throwaway Z300:Z304 helper cells - all the
# "real" destination ranges already received their formatting above.
$ws.Range("Z300:Z304").Clear()

$excel.ActiveWorkbook.Save()
